# Modify growth rate for FCV
#
# The "Max Growth" table (rows 6-17) historically carried the RHS value for
# each UC_CAP constraint in one of two columns: L (projected-2018 base,
# referencing column C of the helper table below) or M (projected-2030,
# referencing column E). The FCV rows (CAR/LGT/MGT/HGT "MaxGrowthFCV") were
# the only ones using the M/"2030" variant.
#
# This edit switches the four FCV rows onto the same 2018-base figure (L,
# referencing C) that every other row already uses, tightens their growth
# multiplier in column J from 1.3 to 1.2, and then drops the now entirely
# unused M column (and its "UC_RHSRTS~2030" header string) so N/O shift
# left into M/N. The FCV rows also get a thin bottom border across B:N to
# set them apart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> the row of the helper "New sales" table (43:46) that its RHS
# formula divides by 1000 (HGT=43, MGT=44, LGV/LGT=45, PC/CAR=46).
$fcvRows = @{8 = 46; 11 = 45; 14 = 44; 17 = 43}

# Step 1: delete column M outright. Every non-FCV row already has M empty,
# so this just removes the dead "UC_RHSRTS~2030" column and shifts the
# "UC_RHSRTS~0" / "UC_Desc" columns (N, O) one to the left (into M, N).
# For the FCV rows this also deletes their old "=-E../1000" formula, which
# step 2 replaces in column L (which was empty for those rows).
$ws.Columns("M").Delete()

foreach ($r in $fcvRows.Keys) {
    $srcRow = $fcvRows[$r]

    # Growth multiplier 1.3 -> 1.2.
    $ws.Range("J$r").Value = 1.2

    # Recreate the RHS formula in L (was empty), now based on the 2018
    # base (column C) instead of the 2030 projection (column E).
    $ws.Range("L$r").Formula = "=-C$srcRow/1000"
    $ws.Range("L$r").NumberFormat = "0.000"

    # Thin black bottom border across the whole row (B:N) to mark it out.
    $rowRng = $ws.Range("B$r`:N$r")
    $rowRng.Borders.Item(9).Color = 0
    $rowRng.Borders.Item(9).LineStyle = 1
}

# Restore the selection Excel had when the file was last saved.
$ws.Range("J21").Select()
